$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.226.53'
$ws.Range('E2').Value = '  -2.06%  '
$ws.Range('D3').Value = '3.482.48'
$ws.Range('E3').Value = '  -1.13%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.14'
$ws.Range('E5').Value = '  -2.93%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.14'
$ws.Range('E6').Value = '  -3.58%  '
$ws.Range('D7').Value = '3.485.06'
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.489'
$ws.Range('E9').Value = '  -4.06%  '
$ws.Range('E10').Value = '  -5.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.24'
$ws.Range('E11').Value = '  -6.44%  '
$ws.Range('E12').Value = '  -6.64%  '
$ws.Range('D13').Value = '4.070.54'
$ws.Range('E13').Value = '  -1.01%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000183'
$ws.Range('E14').Value = '  -6.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '26.64'
$ws.Range('E15').Value = '  -7.10%  '
$ws.Range('D16').Value = '3.477.87'
$ws.Range('E16').Value = '  -1.27%  '
$ws.Range('E17').Value = '  -1.19%  '
$ws.Range('D18').Value = '65.061.43'
$ws.Range('E18').Value = '  -2.06%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.73'
$ws.Range('E19').Value = '  -9.50%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.79'
$ws.Range('E20').Value = '  -5.76%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.89'
$ws.Range('E21').Value = '  -5.04%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '388.63'
$ws.Range('E22').Value = '  -8.12%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.555'
$ws.Range('E23').Value = '  -5.50%  '
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '72.58'
$ws.Range('E25').Value = '  -5.78%  '
$ws.Range('D26').Value = '3.619.66'
$ws.Range('E26').Value = '  -1.22%  '
$ws.Range('E27').Value = '  +0.60%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0000110'
$ws.Range('E28').Value = '  -3.46%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '7.44'
$ws.Range('E30').Value = '  -6.04%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '8.24'
$ws.Range('E31').Value = '  -7.85%  '
$ws.Range('E32').Value = '  -9.89%  '
$ws.Range('D33').Value = '3.494.92'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -6.74%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '23.05'
$ws.Range('E36').Value = '  -4.60%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '172.07'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -8.80%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.88'
$ws.Range('E39').Value = '  -9.16%  '
$ws.Range('E40').Value = '  -9.69%  '
$ws.Range('E41').Value = '  -8.72%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0778'
$ws.Range('E42').Value = '  -4.32%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.813'
$ws.Range('E43').Value = '  -4.49%  '
$ws.Range('E44').Value = '  -0.07%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '42.38'
$ws.Range('E45').Value = '  -6.82%  '
$ws.Range('E46').Value = '  -12.70%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '24.42'
$ws.Range('E47').Value = '  +6.74%  '
$ws.Range('E48').Value = '  -8.51%  '
$ws.Range('E49').Value = '  +1.94%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.67'
$ws.Range('E50').Value = '  -5.74%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.09'
$ws.Range('E51').Value = '  -11.17%  '
